$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. First paragraph: "This is a Microsoft word document." gains a
#    trailing two-space run, followed by three new red-colored runs
#    that read "(This is a change – Version for main branch)".
# ------------------------------------------------------------------
$enDash = [char]0x2013

$findRng = $d.Content
$findRng.Find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRng.InsertAfter("  ")
$findRng.Collapse(0)

$r1 = $findRng.Duplicate
$r1.InsertAfter("(This is a change " + $enDash + " Ve")
$r1.Font.Color = 255
$r1.Collapse(0)

$r2 = $r1.Duplicate
$r2.InsertAfter("rsion for main branch")
$r2.Font.Color = 255
$r2.Collapse(0)

$r3 = $r2.Duplicate
$r3.InsertAfter(")")
$r3.Font.Color = 255
$r3.Collapse(0)

# ------------------------------------------------------------------
# 2. Remove the trailing "ank God almighty, we are free at last."
#    paragraph (the one right after "Shall be lifted—nevermore!").
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Delete()
